$d = $word.ActiveDocument

# Locate the "ЛК № 12 ..." heading paragraph robustly (by content) rather than by a
# hard-coded index.
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*ЛК № 12*") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find the 'ЛК № 12' paragraph"
}

# The heading paragraph is immediately followed by a run of empty "spacer"
# paragraphs, and then a paragraph whose only content is a lone ".".
# Find that lone-dot paragraph after the heading.
$dotIndex = -1
for ($i = $headingIndex + 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq ".`r") {
        $dotIndex = $i
        break
    }
}

if ($dotIndex -eq -1) {
    throw "Could not find the lone '.' paragraph after the heading"
}

# 1) Restore the hanging indent on the lone-dot paragraph (it previously had
#    none; matching its siblings it should have left=2410 twips / hanging=2268
#    twips, i.e. 120.5pt / -113.4pt in the object model).
$dotPara = $d.Paragraphs.Item($dotIndex)
$dotPara.LeftIndent = 120.5
$dotPara.FirstLineIndent = -113.4

# 2) Delete the empty spacer paragraphs that sit between the heading and the
#    lone-dot paragraph (everything strictly between the two).
if ($dotIndex - 1 -ge $headingIndex + 1) {
    $spacerStart = $d.Paragraphs.Item($headingIndex + 1).Range.Start
    $spacerEnd = $d.Paragraphs.Item($dotIndex - 1).Range.End
    $spacerRange = $d.Range($spacerStart, $spacerEnd)
    $spacerRange.Delete()
}

# 3) Clear the heading paragraph's text (runs), leaving it as an empty
#    paragraph while keeping its own paragraph mark/formatting.
$headingPara = $d.Paragraphs.Item($headingIndex)
$headingTextRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End - 1)
$headingTextRange.Text = ""
